$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2239496666666667
$ws.Range("H2").Value = 0.671849
$ws.Range("I2").Value = 0.4264743968982249
$ws.Range("J2").Value = 0.4264743968982249
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01518366666666667
$ws.Range("N2").Value = 0.045551
$ws.Range("O2").Value = 0.02370341769240456
$ws.Range("P2").Value = 0.02370341769240456
$ws.Range("Q2").Value = 0.003400377088777778
$ws.Range("R2").Value = 0.030603393799
$ws.Range("S2").Value = 0.01010890076479495
$ws.Range("T2").Value = 0.01010890076479495
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2239496666666667
$ws.Range("H3").Value = 0.671849
$ws.Range("I3").Value = 0.4264743968982249
$ws.Range("J3").Value = 0.4264743968982249
$ws.Range("O3").Value = 0.8698427334878488
$ws.Range("P3").Value = 0.8698427334878488
$ws.Range("Q3").Value = 0.1247834105686667
$ws.Range("R3").Value = 1.123050695118
$ws.Range("S3").Value = 0.3709656551605337
$ws.Range("T3").Value = 0.3709656551605336
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2239496666666667
$ws.Range("H4").Value = 0.671849
$ws.Range("I4").Value = 0.4264743968982249
$ws.Range("J4").Value = 0.4264743968982249
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.068191
$ws.Range("N4").Value = 0.204573
$ws.Range("O4").Value = 0.1064538488197466
$ws.Range("P4").Value = 0.1064538488197466
$ws.Range("Q4").Value = 0.01527135171966667
$ws.Range("R4").Value = 0.137442165477
$ws.Range("S4").Value = 0.04539984097289625
$ws.Range("T4").Value = 0.04539984097289625
$ws.Range("I5").Value = 0.4001470143891285
$ws.Range("J5").Value = 0.4001470143891285
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01518366666666667
$ws.Range("N5").Value = 0.045551
$ws.Range("O5").Value = 0.02370341769240456
$ws.Range("P5").Value = 0.02370341769240456
$ws.Range("Q5").Value = 0.003190462897111111
$ws.Range("R5").Value = 0.028714166074
$ws.Range("S5").Value = 0.009484851820434131
$ws.Range("T5").Value = 0.009484851820434131
$ws.Range("I6").Value = 0.4001470143891285
$ws.Range("J6").Value = 0.4001470143891285
$ws.Range("O6").Value = 0.8698427334878488
$ws.Range("P6").Value = 0.8698427334878488
$ws.Range("S6").Value = 0.3480649727932411
$ws.Range("T6").Value = 0.3480649727932411
$ws.Range("I7").Value = 0.4001470143891285
$ws.Range("J7").Value = 0.4001470143891285
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.068191
$ws.Range("N7").Value = 0.204573
$ws.Range("O7").Value = 0.1064538488197466
$ws.Range("P7").Value = 0.1064538488197466
$ws.Range("Q7").Value = 0.01432861114466667
$ws.Range("R7").Value = 0.128957500302
$ws.Range("S7").Value = 0.04259718977545326
$ws.Range("T7").Value = 0.04259718977545326
$ws.Range("G8").Value = 0.09104433333333334
$ws.Range("H8").Value = 0.273133
$ws.Range("I8").Value = 0.1733785887126465
$ws.Range("J8").Value = 0.1733785887126465
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01518366666666667
$ws.Range("N8").Value = 0.045551
$ws.Range("O8").Value = 0.02370341769240456
$ws.Range("P8").Value = 0.02370341769240456
$ws.Range("Q8").Value = 0.001382386809222222
$ws.Range("R8").Value = 0.012441481283
$ws.Range("S8").Value = 0.00410966510717548
$ws.Range("T8").Value = 0.00410966510717548
$ws.Range("G9").Value = 0.09104433333333334
$ws.Range("H9").Value = 0.273133
$ws.Range("I9").Value = 0.1733785887126465
$ws.Range("J9").Value = 0.1733785887126465
$ws.Range("O9").Value = 0.8698427334878488
$ws.Range("P9").Value = 0.8698427334878488
$ws.Range("Q9").Value = 0.05072935626733333
$ws.Range("R9").Value = 0.456564206406
$ws.Range("S9").Value = 0.150812105534074
$ws.Range("T9").Value = 0.150812105534074
$ws.Range("G10").Value = 0.09104433333333334
$ws.Range("H10").Value = 0.273133
$ws.Range("I10").Value = 0.1733785887126465
$ws.Range("J10").Value = 0.1733785887126465
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.068191
$ws.Range("N10").Value = 0.204573
$ws.Range("O10").Value = 0.1064538488197466
$ws.Range("P10").Value = 0.1064538488197466
$ws.Range("Q10").Value = 0.006208404134333334
$ws.Range("R10").Value = 0.055875637209
$ws.Range("S10").Value = 0.0184568180713971
$ws.Range("T10").Value = 0.0184568180713971
